# Weekly update: a new daily price record is inserted as row 130 (pushing
# all subsequent rows down by one), and it is populated with a new date and
# a new "Volumen" (J) value while carrying the rest of the template values
# forward from the row it was inserted above.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 130; existing rows 130:225 shift to 131:226.
$ws.Rows(130).Insert()

# Populate the newly inserted row 130 with the new data point.
$ws.Range("A130").Value = 8
$ws.Range("B130").Value = "Terminal La Palmera de La Serena"
$ws.Range("C130").Value = "Coquimbo"
$ws.Range("D130").Value = "2023-05-16"
$ws.Range("E130").Value = 4
$ws.Range("F130").Value = 100112044
$ws.Range("G130").Value = "Perejil"
$ws.Range("H130").Value = "Sin especificar"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 2200
$ws.Range("K130").Value = 2000
$ws.Range("L130").Value = 2500
$ws.Range("M130").Value = 2250
$ws.Range("N130").Value = "$/atado 1 a 1,5 kilos"
$ws.Range("O130").Value = "Provincia del Elquí"
$ws.Range("P130").Value = 1500
$ws.Range("Q130").Value = 1.5
$ws.Range("R130").Value = "Hortaliza"
